$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data for rows 19-23: MSSV, Ho va ten, and 7 "Phan cong" percentage columns (E..K)
$data = @(
    @{ Row = 19; MSSV = "0712152"; Name = "Lê Long Hồ" },
    @{ Row = 20; MSSV = "0712163"; Name = "Võ Minh Hiển" },
    @{ Row = 21; MSSV = "0712174"; Name = "Nguyễn Văn Hiếu" },
    @{ Row = 22; MSSV = "0712178"; Name = "Nguyễn Ngọc Hoà" },
    @{ Row = 23; MSSV = "0712190"; Name = "Lê Gia Quốc Huy" }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 3).Value = $entry.MSSV
    $ws.Cells.Item($r, 4).Value = $entry.Name

    # Columns E (5) through K (11): 100% completion, formatted as percentage
    for ($c = 5; $c -le 11; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Value = 1
        $cell.NumberFormat = "0%"
    }
}

$ws.Range("K19").Select()
